# Natmi following Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Cells.Item(2,1).Value = "ECs"
    $ws.Cells.Item(2,2).Value = "Cthrc1"
    $ws.Cells.Item(2,3).Value = "Fzd6"
    $ws.Cells.Item(2,4).Value = "ECs"
    $ws.Cells.Item(2,5).Value = 2
    $ws.Cells.Item(2,6).Value = 0.6666666666666666
    $ws.Cells.Item(2,7).Value = 0.5226613333333333
    $ws.Cells.Item(2,8).Value = 1.567984
    $ws.Cells.Item(2,9).Value = 0.004118772370031606
    $ws.Cells.Item(2,10).Value = 0.004118772370031606
    $ws.Cells.Item(2,11).Value = 3
    $ws.Cells.Item(2,12).Value = 1
    $ws.Cells.Item(2,13).Value = 12.415956
    $ws.Cells.Item(2,14).Value = 37.247868
    $ws.Cells.Item(2,15).Value = 0.8598042313376485
    $ws.Cells.Item(2,16).Value = 0.8598042313376485
    $ws.Cells.Item(2,17).Value = 6.489340117567999
    $ws.Cells.Item(2,18).Value = 58.404061058112
    $ws.Cells.Item(2,19).Value = 0.003541337911669769
    $ws.Cells.Item(2,20).Value = 0.003541337911669769

    # Row 3
    $ws.Cells.Item(3,1).Value = "ECs"
    $ws.Cells.Item(3,2).Value = "Cthrc1"
    $ws.Cells.Item(3,3).Value = "Fzd6"
    $ws.Cells.Item(3,4).Value = "FAPs"
    $ws.Cells.Item(3,5).Value = 2
    $ws.Cells.Item(3,6).Value = 0.6666666666666666
    $ws.Cells.Item(3,7).Value = 0.5226613333333333
    $ws.Cells.Item(3,8).Value = 1.567984
    $ws.Cells.Item(3,9).Value = 0.004118772370031606
    $ws.Cells.Item(3,10).Value = 0.004118772370031606
    $ws.Cells.Item(3,11).Value = 3
    $ws.Cells.Item(3,12).Value = 1
    $ws.Cells.Item(3,13).Value = 1.745879666666666
    $ws.Cells.Item(3,14).Value = 5.237639
    $ws.Cells.Item(3,15).Value = 0.1209020654395331
    $ws.Cells.Item(3,16).Value = 0.1209020654395331
    $ws.Cells.Item(3,17).Value = 0.9125037944195554
    $ws.Cells.Item(3,18).Value = 8.212534149775999
    $ws.Cells.Item(3,19).Value = 0.0004979680866121019
    $ws.Cells.Item(3,20).Value = 0.0004979680866121019

    # Row 4
    $ws.Cells.Item(4,1).Value = "ECs"
    $ws.Cells.Item(4,2).Value = "Cthrc1"
    $ws.Cells.Item(4,3).Value = "Fzd6"
    $ws.Cells.Item(4,4).Value = "sCs"
    $ws.Cells.Item(4,5).Value = 2
    $ws.Cells.Item(4,6).Value = 0.6666666666666666
    $ws.Cells.Item(4,7).Value = 0.5226613333333333
    $ws.Cells.Item(4,8).Value = 1.567984
    $ws.Cells.Item(4,9).Value = 0.004118772370031606
    $ws.Cells.Item(4,10).Value = 0.004118772370031606
    $ws.Cells.Item(4,11).Value = 2
    $ws.Cells.Item(4,12).Value = 0.6666666666666666
    $ws.Cells.Item(4,13).Value = 0.2786096666666666
    $ws.Cells.Item(4,14).Value = 0.8358289999999999
    $ws.Cells.Item(4,15).Value = 0.01929370322281843
    $ws.Cells.Item(4,16).Value = 0.01929370322281843
    $ws.Cells.Item(4,17).Value = 0.1456184998595556
    $ws.Cells.Item(4,18).Value = 1.310566498736
    $ws.Cells.Item(4,19).Value = 0.00007946637174973428
    $ws.Cells.Item(4,20).Value = 0.0000794663717497343

    # Row 5
    $ws.Cells.Item(5,1).Value = "FAPs"
    $ws.Cells.Item(5,2).Value = "Cthrc1"
    $ws.Cells.Item(5,3).Value = "Fzd6"
    $ws.Cells.Item(5,4).Value = "ECs"
    $ws.Cells.Item(5,5).Value = 3
    $ws.Cells.Item(5,6).Value = 1
    $ws.Cells.Item(5,7).Value = 125.8872733333333
    $ws.Cells.Item(5,8).Value = 377.66182
    $ws.Cells.Item(5,9).Value = 0.9920401416289004
    $ws.Cells.Item(5,10).Value = 0.9920401416289004
    $ws.Cells.Item(5,11).Value = 3
    $ws.Cells.Item(5,12).Value = 1
    $ws.Cells.Item(5,13).Value = 12.415956
    $ws.Cells.Item(5,14).Value = 37.247868
    $ws.Cells.Item(5,15).Value = 0.8598042313376485
    $ws.Cells.Item(5,16).Value = 0.8598042313376485
    $ws.Cells.Item(5,17).Value = 1563.01084666664
    $ws.Cells.Item(5,18).Value = 14067.09761999976
    $ws.Cells.Item(5,19).Value = 0.8529603114293287
    $ws.Cells.Item(5,20).Value = 0.8529603114293287

    # Row 6
    $ws.Cells.Item(6,1).Value = "FAPs"
    $ws.Cells.Item(6,2).Value = "Cthrc1"
    $ws.Cells.Item(6,3).Value = "Fzd6"
    $ws.Cells.Item(6,4).Value = "FAPs"
    $ws.Cells.Item(6,5).Value = 3
    $ws.Cells.Item(6,6).Value = 1
    $ws.Cells.Item(6,7).Value = 125.8872733333333
    $ws.Cells.Item(6,8).Value = 377.66182
    $ws.Cells.Item(6,9).Value = 0.9920401416289004
    $ws.Cells.Item(6,10).Value = 0.9920401416289004
    $ws.Cells.Item(6,11).Value = 3
    $ws.Cells.Item(6,12).Value = 1
    $ws.Cells.Item(6,13).Value = 1.745879666666666
    $ws.Cells.Item(6,14).Value = 5.237639
    $ws.Cells.Item(6,15).Value = 0.1209020654395331
    $ws.Cells.Item(6,16).Value = 0.1209020654395331
    $ws.Cells.Item(6,17).Value = 219.7840308047755
    $ws.Cells.Item(6,18).Value = 1978.05627724298
    $ws.Cells.Item(6,19).Value = 0.119939702121861
    $ws.Cells.Item(6,20).Value = 0.119939702121861

    # Row 7
    $ws.Cells.Item(7,1).Value = "FAPs"
    $ws.Cells.Item(7,2).Value = "Cthrc1"
    $ws.Cells.Item(7,3).Value = "Fzd6"
    $ws.Cells.Item(7,4).Value = "sCs"
    $ws.Cells.Item(7,5).Value = 3
    $ws.Cells.Item(7,6).Value = 1
    $ws.Cells.Item(7,7).Value = 125.8872733333333
    $ws.Cells.Item(7,8).Value = 377.66182
    $ws.Cells.Item(7,9).Value = 0.9920401416289004
    $ws.Cells.Item(7,10).Value = 0.9920401416289004
    $ws.Cells.Item(7,11).Value = 2
    $ws.Cells.Item(7,12).Value = 0.6666666666666666
    $ws.Cells.Item(7,13).Value = 0.2786096666666666
    $ws.Cells.Item(7,14).Value = 0.8358289999999999
    $ws.Cells.Item(7,15).Value = 0.01929370322281843
    $ws.Cells.Item(7,16).Value = 0.01929370322281843
    $ws.Cells.Item(7,17).Value = 35.07341126097555
    $ws.Cells.Item(7,18).Value = 315.66070134878
    $ws.Cells.Item(7,19).Value = 0.01914012807771076
    $ws.Cells.Item(7,20).Value = 0.01914012807771077

    # Row 8
    $ws.Cells.Item(8,1).Value = "sCs"
    $ws.Cells.Item(8,2).Value = "Cthrc1"
    $ws.Cells.Item(8,3).Value = "Fzd6"
    $ws.Cells.Item(8,4).Value = "ECs"
    $ws.Cells.Item(8,5).Value = 2
    $ws.Cells.Item(8,6).Value = 0.6666666666666666
    $ws.Cells.Item(8,7).Value = 0.4874236666666666
    $ws.Cells.Item(8,8).Value = 1.462271
    $ws.Cells.Item(8,9).Value = 0.003841086001067923
    $ws.Cells.Item(8,10).Value = 0.003841086001067923
    $ws.Cells.Item(8,11).Value = 3
    $ws.Cells.Item(8,12).Value = 1
    $ws.Cells.Item(8,13).Value = 12.415956
    $ws.Cells.Item(8,14).Value = 37.247868
    $ws.Cells.Item(8,15).Value = 0.8598042313376485
    $ws.Cells.Item(8,16).Value = 0.8598042313376485
    $ws.Cells.Item(8,17).Value = 6.051830798691999
    $ws.Cells.Item(8,18).Value = 54.46647718822799
    $ws.Cells.Item(8,19).Value = 0.003302581996650008
    $ws.Cells.Item(8,20).Value = 0.003302581996650007

    # Row 9
    $ws.Cells.Item(9,1).Value = "sCs"
    $ws.Cells.Item(9,2).Value = "Cthrc1"
    $ws.Cells.Item(9,3).Value = "Fzd6"
    $ws.Cells.Item(9,4).Value = "FAPs"
    $ws.Cells.Item(9,5).Value = 2
    $ws.Cells.Item(9,6).Value = 0.6666666666666666
    $ws.Cells.Item(9,7).Value = 0.4874236666666666
    $ws.Cells.Item(9,8).Value = 1.462271
    $ws.Cells.Item(9,9).Value = 0.003841086001067923
    $ws.Cells.Item(9,10).Value = 0.003841086001067923
    $ws.Cells.Item(9,11).Value = 3
    $ws.Cells.Item(9,12).Value = 1
    $ws.Cells.Item(9,13).Value = 1.745879666666666
    $ws.Cells.Item(9,14).Value = 5.237639
    $ws.Cells.Item(9,15).Value = 0.1209020654395331
    $ws.Cells.Item(9,16).Value = 0.1209020654395331
    $ws.Cells.Item(9,17).Value = 0.8509830686854444
    $ws.Cells.Item(9,18).Value = 7.658847618168999
    $ws.Cells.Item(9,19).Value = 0.0004643952310599884
    $ws.Cells.Item(9,20).Value = 0.0004643952310599884

    # Row 10
    $ws.Cells.Item(10,1).Value = "sCs"
    $ws.Cells.Item(10,2).Value = "Cthrc1"
    $ws.Cells.Item(10,3).Value = "Fzd6"
    $ws.Cells.Item(10,4).Value = "sCs"
    $ws.Cells.Item(10,5).Value = 2
    $ws.Cells.Item(10,6).Value = 0.6666666666666666
    $ws.Cells.Item(10,7).Value = 0.4874236666666666
    $ws.Cells.Item(10,8).Value = 1.462271
    $ws.Cells.Item(10,9).Value = 0.003841086001067923
    $ws.Cells.Item(10,10).Value = 0.003841086001067923
    $ws.Cells.Item(10,11).Value = 2
    $ws.Cells.Item(10,12).Value = 0.6666666666666666
    $ws.Cells.Item(10,13).Value = 0.2786096666666666
    $ws.Cells.Item(10,14).Value = 0.8358289999999999
    $ws.Cells.Item(10,15).Value = 0.01929370322281843
    $ws.Cells.Item(10,16).Value = 0.01929370322281843
    $ws.Cells.Item(10,17).Value = 0.1358009452954444
    $ws.Cells.Item(10,18).Value = 1.222208507659
    $ws.Cells.Item(10,19).Value = 0.00007410877335792692
    $ws.Cells.Item(10,20).Value = 0.00007410877335792694
